$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.681.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.812.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.560'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.20%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.57%  '

$ws.Range("E11").Value = '  -0.83%  '

$ws.Range("E12").Value = '  +2.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.260.32'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.973'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.812.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.742.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.05%  '

$ws.Range("E19").Value = '  +10.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0971'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.55%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.90%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.71%  '

$ws.Range("E31").Value = '  +2.63%  '

$ws.Range("E32").Value = '  +1.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.73%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0899'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.91%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0453'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.40%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.08%  '

$ws.Range("E40").Value = '  +2.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.38%  '

$ws.Range("E42").Value = '  +1.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '121.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.80%  '

$ws.Range("E44").Value = '  +1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.25%  '

$ws.Range("E47").Value = '  +8.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.153.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.985'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.225'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +17.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0321'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.74%  '
